# Re-sort the sheet tabs: put "总计" (the summary sheet) before "2022-Q1"
# (the detail sheet), i.e. move it to be the first sheet in the workbook.
$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$firstSheet = $wb.Worksheets.Item(1)

$summarySheet.Move($firstSheet)
